$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" ---
$wsMeta = $wb.Worksheets.Item("Metadata")

# Version bump
$wsMeta.Range("B3").Value = "6.0.0"

# Date update
$wsMeta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value now populated
$wsMeta.Range("B9").Value = "Alvearie Team"

# Row 10 was "Contact" / "No display for ContactDetail" -> becomes "Jurisdiction" / "United States of America"
$wsMeta.Range("A10").Value = "Jurisdiction"
$wsMeta.Range("B10").Value = "United States of America"

# Row 11 was a duplicate "Contact" / "No display for ContactDetail" row -> remove it entirely,
# shifting the remaining rows (old 12-21) up to become 11-20.
$wsMeta.Rows.Item(11).Delete()

# --- Sheet "Elements" ---
$wsElements = $wb.Worksheets.Item("Elements")

# Row 2 (the Extension element itself): Short/Definition now reflect this specific extension
$wsElements.Range("K2").Value = "Reinsurance Met Indicator"
$wsElements.Range("L2").Value = "Indicates whether the reinsurance amount was met"
